$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new "Sample-type" column before F; old F (Impact on Baseline) shifts to G ---
$ws.Columns("F:F").Insert()

# --- Header row ---
$ws.Range("F1").Value = "Sample-type"

# --- Update existing scenario descriptions to sputum-based wording ---
$ws.Range("B3").Value = "Individuals tested through a sputum-based GeneXpert. Testing can occur onsite but is predominantly centralised and perfromed offsite. Test results are provided at subsequent visits. "
$ws.Range("B4").Value = "Everyone receives a decentralised sputum-based GeneXpert. Results provided at subsequent visits"
$ws.Range("B5").Value = "Everyone receives a decentralised sputum-based GeneXert and gets result in same visit"

# --- Fill new Sample-type column for existing rows 2-5 ---
$ws.Range("F2").Value = "Sputum"
$ws.Range("F2").WrapText = $true
$ws.Range("F3").Value = "Sputum"
$ws.Range("F3").WrapText = $true
$ws.Range("F4").Value = "Sputum"
$ws.Range("F4").WrapText = $true
$ws.Range("F5").Value = "Sputum"
$ws.Range("F5").WrapText = $true

# --- Insert 9 new rows (6-14) below row 5 for new scenarios + blank spacer rows ---
$ws.Rows("6:14").Insert()

# --- Remove stray formatted-but-empty cells that the row insert carried into the blank spacer rows ---
$ws.Range("B9:D12").Clear()

# --- Row 6 ---
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "Individuals tested through a swab-based GeneXpert. Testing can occur onsite but is predominantly centralised and perfromed offsite. Test results are provided at subsequent visits. "
$ws.Range("B6").WrapText = $true
$ws.Range("C6").Value = "GeneXpert"
$ws.Range("D6").Value = "Mix of onsite and off site (centralised)"
$ws.Range("D6").WrapText = $true
$ws.Range("E6").Value = "Subsequent visit"
$ws.Range("F6").Value = "Non-sputum (oral swab)"
$ws.Range("F6").WrapText = $true
$ws.Range("G6").Value = "Testing occurs fully with GeneXpert, as opposed to a proportion with smear microscopy, hence a greater number of correct results will be generated. Oral swab will allow for all individuals to provide a sample (particulalry HIV+)  and increase ease and likelihood of testing being offered"
$ws.Range("G6").WrapText = $true

# --- Row 7 ---
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "Everyone receives a decentralised swab-based GeneXpert. Results provided at subsequent visits"
$ws.Range("B7").WrapText = $true
$ws.Range("C7").Value = "GeneXpert"
$ws.Range("D7").Value = "Onsite (decentralised)"
$ws.Range("D7").WrapText = $true
$ws.Range("E7").Value = "Subsequent visit"
$ws.Range("F7").Value = "Non-sputum (oral swab)"
$ws.Range("F7").WrapText = $true
$ws.Range("G7").Value = "A greater number of correct results are generated. Further, all sites have testing capapcity, hence fewer people are lost from care  as a result of needing to attend another facility. Oral swab will allow for all individuals to provide a sample (particulalry HIV+) and increase ease and likelihood of testing being offered"
$ws.Range("G7").WrapText = $true

# --- Row 8 ---
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "Everyone receives a decentralised swab-based GeneXert and gets result in same visit"
$ws.Range("B8").WrapText = $true
$ws.Range("C8").Value = "GeneXpert"
$ws.Range("D8").Value = "Onsite (decentralised)"
$ws.Range("D8").WrapText = $true
$ws.Range("E8").Value = "At testing visit"
$ws.Range("F8").Value = "Non-sputum (oral swab)"
$ws.Range("F8").WrapText = $true
$ws.Range("G8").Value = "A greater number of correct results are generated. All sites have testing capapcity, hence fewer people are lost from care as a result of needing to attend another facility. All results are availble at the time of testing and a subsequent visit to collect results is not required. Oral swab will allow for all individuals to provide a sample (particulalry HIV+)  and increase ease and likelihood of testing being offered"
$ws.Range("G8").WrapText = $true

# --- Blank spacer rows 9-12: just wrap-formatted empty cells in F/G ---
$ws.Range("F9").WrapText = $true
$ws.Range("G9").WrapText = $true
$ws.Range("F10").WrapText = $true
$ws.Range("G10").WrapText = $true
$ws.Range("F11").WrapText = $true
$ws.Range("G11").WrapText = $true
$ws.Range("F12").WrapText = $true
$ws.Range("G12").WrapText = $true

# --- Row 13 ---
$ws.Range("A13").Value = 4
$ws.Range("B13").Value = "Everone receives a low-compexity sputum-based test (98% specificity, 91% sensitivity). Similar/Equivalent to GeneXert use case in Scenario 1"
$ws.Range("B13").WrapText = $true
$ws.Range("C13").Value = "TPP test - sputum, low-complexity Assay"
$ws.Range("D13").Value = "Mix of onsite and off site (centralised)"
$ws.Range("D13").WrapText = $true
$ws.Range("E13").Value = "Subsequent visit"
$ws.Range("F13").Value = "Sputum"
$ws.Range("F13").WrapText = $true
$ws.Range("G13").Value = "Testing occurs fully with low-complexity assay, as opposed to a proportion with smear microscopy, hence a greater number of correct results will be generated."
$ws.Range("G13").WrapText = $true

# --- Row 14 ---
$ws.Range("A14").Value = 4
$ws.Range("B14").Value = "Everyone receives a low-compexity non-sputum-based test (98% specificity, 77% sensitivity)."
$ws.Range("B14").WrapText = $true
$ws.Range("C14").Value = "TPP test - non-sputum, low-complexity Assay"
$ws.Range("D14").Value = "Mix of onsite and off site (centralised)"
$ws.Range("D14").WrapText = $true
$ws.Range("E14").Value = "Subsequent visit"
$ws.Range("F14").Value = "Non-sputum"
$ws.Range("F14").WrapText = $true
$ws.Range("G14").Value = "Testing occurs fully with low-complexity assay, as opposed to a proportion with smear microscopy, hence a greater number of correct results will be generated. Non-sputum samples increase the numer of individuals with HIV who can provide a sample. "
$ws.Range("G14").WrapText = $true

# --- Row heights (Excel normally auto-sizes these from wrapped content; set explicitly) ---
$ws.Rows("1:1").RowHeight = 17
$ws.Rows("6:6").RowHeight = 187
$ws.Rows("7:7").RowHeight = 187
$ws.Rows("8:8").RowHeight = 238
$ws.Rows("13:13").RowHeight = 102
$ws.Rows("14:14").RowHeight = 153

# --- Sheet view / selection adjustments ---
$ws.Range("A7").Select()
$excel.ActiveWindow.ScrollRow = 7
$ws.Range("J9").Select()

